# Adding process edit address page
# - Rename the existing LOCAL_URL= env entry to LOCAL_URL_AD= (admin URL)
# - Insert a new LOCAL_URL_CL= entry (client URL) right below it
#
# This shifts every row at/after row 8 down by one, which matches all of
# the row renumbering seen further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 8 (pushes old rows 8.. down by one)
$ws.Rows(8).Insert()

# Row 7: rename LOCAL_URL= -> LOCAL_URL_AD= (value stays the same URL)
$r7 = $ws.Range("A7")
$r7.Value = "LOCAL_URL_AD=http://localhost:3005"
$r7.Characters(14, 22).Font.Color = 255
$r7.Characters(14, 22).Font.Name = "Calibri"
$r7.Characters(14, 22).Font.Size = 11

# Row 8 (new): LOCAL_URL_CL= for the client URL
$r8 = $ws.Range("A8")
$r8.Value = "LOCAL_URL_CL=http://localhost:3005"
$r8.Characters(14, 22).Font.Color = 255
$r8.Characters(14, 22).Font.Name = "Calibri"
$r8.Characters(14, 22).Font.Size = 11

# Match the saved selection state from the authored workbook
$ws.Range("K23").Select()
